$d = $word.ActiveDocument

# Each of the three "<id>...</id>" blocks was split across three runs:
#   <id>   (Courier New, color 7f6000)
#   p130v_N (color 000000)
#   </id>  (Courier New, color 7f6000)
# Merge them back into a single run (keeping the first run's formatting)
# by finding the full "<id>p130v_N</id>" text and replacing it with itself.

$ids = @("p130v_1", "p130v_2", "p130v_3")

foreach ($id in $ids) {
    $needle = "<id>" + $id + "</id>"
    $rng = $d.Content
    $found = $rng.Find.Execute($needle, $false, $false, $false, $false, $false, $true, 1, $false, $needle, 2)
}
